$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.509.86"
$ws.Range("E2").Value = "  -0.52%  "

$ws.Range("D3").Value = "2.945.18"
$ws.Range("E3").Value = "  -2.11%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'570.85"
$ws.Range("E5").Value = "  -2.02%  "

$ws.Range("D6").Value = "'162.58"
$ws.Range("E6").Value = "  +0.47%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").Value = "'0.516"
$ws.Range("E8").Value = "  -0.59%  "

$ws.Range("D9").Value = "2.942.10"
$ws.Range("E9").Value = "  -2.12%  "

$ws.Range("E10").Value = "  -0.99%  "

$ws.Range("E11").Value = "  -4.06%  "

$ws.Range("D12").Value = "'0.464"
$ws.Range("E12").Value = "  +0.98%  "

$ws.Range("D13").Value = "'0.0000243"
$ws.Range("E13").Value = "  -3.04%  "

$ws.Range("D14").Value = "'34.78"
$ws.Range("E14").Value = "  +0.52%  "

$ws.Range("E15").Value = "  -0.55%  "

$ws.Range("D16").Value = "65.559.80"
$ws.Range("E16").Value = "  -0.45%  "

$ws.Range("D17").Value = "3.437.77"
$ws.Range("E17").Value = "  -1.97%  "

$ws.Range("D18").Value = "'7.03"
$ws.Range("E18").Value = "  +0.76%  "

$ws.Range("D19").Value = "2.950.22"
$ws.Range("E19").Value = "  -1.87%  "

$ws.Range("D20").Value = "'15.69"
$ws.Range("E20").Value = "  +12.41%  "

$ws.Range("D21").Value = "'444.88"
$ws.Range("E21").Value = "  -2.77%  "

$ws.Range("D22").Value = "'0.694"
$ws.Range("E22").Value = "  +0.75%  "

$ws.Range("D23").Value = "'7.27"
$ws.Range("E23").Value = "  -1.22%  "

$ws.Range("D24").Value = "'82.11"
$ws.Range("E24").Value = "  -0.27%  "

$ws.Range("D25").Value = "'2.26"
$ws.Range("E25").Value = "  -1.50%  "

$ws.Range("D26").Value = "'12.22"
$ws.Range("E26").Value = "  -1.31%  "

$ws.Range("D27").Value = "'10.01"
$ws.Range("E27").Value = "  -6.04%  "

$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("D29").Value = "'2.48"
$ws.Range("E29").Value = "  +5.80%  "

$ws.Range("D30").Value = "'8.05"
$ws.Range("E30").Value = "  -1.14%  "

$ws.Range("E31").Value = "  -0.28%  "

$ws.Range("E32").Value = "  -2.42%  "

$ws.Range("E33").Value = "  +3.69%  "

$ws.Range("D34").Value = "'27.14"
$ws.Range("E34").Value = "  +0.69%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("D36").Value = "'0.973"
$ws.Range("E36").Value = "  -2.22%  "

$ws.Range("D37").Value = "'5.71"
$ws.Range("E37").Value = "  -1.70%  "

$ws.Range("D38").Value = "'46.36"
$ws.Range("E38").Value = "  +5.44%  "

$ws.Range("D39").Value = "'49.06"
$ws.Range("E39").Value = "  -1.33%  "

$ws.Range("E40").Value = "  -1.65%  "

$ws.Range("E41").Value = "  -8.95%  "

$ws.Range("E42").Value = "  -0.42%  "

$ws.Range("E43").Value = "  -4.83%  "

$ws.Range("D44").Value = "'8.50"
$ws.Range("E44").Value = "  +0.65%  "

$ws.Range("D45").Value = "'384.11"
$ws.Range("E45").Value = "  -2.21%  "

$ws.Range("E46").Value = "  -1.31%  "

$ws.Range("D47").Value = "2.676.72"
$ws.Range("E47").Value = "  -4.05%  "

$ws.Range("D48").Value = "'133.59"
$ws.Range("E48").Value = "  -1.04%  "

$ws.Range("D50").Value = "'23.95"
$ws.Range("E50").Value = "  +0.49%  "

$ws.Range("E51").Value = "  +1.10%  "
